# Update workbook for "Add data for 2022-11-21" commit.
# This refreshes the carjacking-by-neighborhood-by-month report from
# "through November 12" to "through November 13", which both renames
# the sheet/header label and adds the newly-tallied incidents to the
# appropriate historical (neighborhood, month) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab and update the column-B header label text to
# reflect the new "as of" date.
$ws.Name = "Through 2022-11-13"
$ws.Range("B1").Value = "November 2022 (through November 13)"

# Per-cell incident-count updates (row = neighborhood, column = month).
$updates = @(
    @{ Cell = "BE2";  Value = 1 }
    @{ Cell = "BP2";  Value = 3 }
    @{ Cell = "AT5";  Value = 2 }
    @{ Cell = "BE5";  Value = 2 }
    @{ Cell = "BP5";  Value = 3 }
    @{ Cell = "BE6";  Value = 1 }
    @{ Cell = "CA6";  Value = 1 }
    @{ Cell = "X9";   Value = 1 }
    @{ Cell = "M16";  Value = 7 }
    @{ Cell = "X16";  Value = 9 }
    @{ Cell = "X17";  Value = 5 }
    @{ Cell = "B18";  Value = 2 }
    @{ Cell = "AI18"; Value = 1 }
    @{ Cell = "BE18"; Value = 4 }
    @{ Cell = "BE21"; Value = 6 }
    @{ Cell = "M24";  Value = 2 }
    @{ Cell = "B32";  Value = 1 }
    @{ Cell = "AT35"; Value = 3 }
    @{ Cell = "BE35"; Value = 1 }
    @{ Cell = "B38";  Value = 1 }
    @{ Cell = "B40";  Value = 2 }
    @{ Cell = "BE40"; Value = 3 }
    @{ Cell = "B56";  Value = 1 }
    @{ Cell = "BE62"; Value = 2 }
    @{ Cell = "BP62"; Value = 1 }
    @{ Cell = "X71";  Value = 1 }
    @{ Cell = "M74";  Value = 2 }
    @{ Cell = "B98";  Value = 1 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
